$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
# Copy formatting from an existing header cell (A1:AC1) so the new header
# cells share the same bold/centered/bordered style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-52: team record Wins/Losses/Ties values
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 59
    $ws.Cells.Item($r, 31).Value = 102
    $ws.Cells.Item($r, 32).Value = 0
}
